$wb = $excel.ActiveWorkbook

# Work on the "DBS" sheet (sheet2) - add a new row describing the
# "findOnlyTran" lookup, mirroring the existing "findByTran" row.
$ws = $wb.Worksheets.Item("DBS")

$ws.Range("B3").Value = "TranNo = "
$ws.Range("A3").Value = "findOnlyTran"
$ws.Range("C3").Value = "FileNo Desc"

# Copy the formatting/style from row 2 so the new row matches.
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A3:C3").PasteSpecial(-4122) | Out-Null

# Make the DBS sheet the active sheet/tab (it was DBD before).
$ws.Activate()
$ws.Range("B9").Select() | Out-Null
